# Update column G (K) values on the active worksheet to reflect
# regenerated save_data (K instead of Strike#, recalculated s_vals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = 0
    4  = 2
    6  = 0
    7  = 1
    8  = 6
    9  = 4
    10 = 1
    11 = 9
    12 = 5
    13 = 7
    14 = 3
    15 = 5
    16 = 6
    17 = 4
    18 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
